$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 129.8
$ws.Range("I4").Value = 136.75
$ws.Range("K4").Value = 136.75
$ws.Range("M4").Value = -22.75
$ws.Range("H12").Value = 1197.5
$ws.Range("J12").Value = 2180.4
$ws.Range("L12").Value = 2180.4
$ws.Range("N12").Value = -2520.4
$ws.Range("H28").Value = 729.3889
$ws.Range("I28").Value = 729.3889
$ws.Range("K28").Value = 729.3889
$ws.Range("M28").Value = -244.3889
$ws.Range("H33").Value = 4167754.5
$ws.Range("I33").Value = 6250132
$ws.Range("J33").Value = 2999.625
$ws.Range("K33").Value = 6250132
$ws.Range("L33").Value = 2999.625
$ws.Range("M33").Value = -6249903
$ws.Range("N33").Value = -3457.625
$ws.Range("H55").Value = 218.33333
$ws.Range("I55").Value = 235.5
$ws.Range("J55").Value = 184
$ws.Range("K55").Value = 235.5
$ws.Range("L55").Value = 184
$ws.Range("M55").Value = -21.5
$ws.Range("N55").Value = -612
$ws.Range("H76").Value = 4715.8335
$ws.Range("H79").Value = 4715.8335
$ws.Range("H80").Value = 874.75
$ws.Range("I80").Value = 668
$ws.Range("K80").Value = 2004
$ws.Range("M80").Value = -1006
$ws.Range("H83").Value = 874.75
$ws.Range("I83").Value = 668
$ws.Range("K83").Value = 6012
$ws.Range("M83").Value = -1020
$ws.Range("H98").Value = 1346.4286
$ws.Range("J98").Value = 1250
$ws.Range("L98").Value = 1250
$ws.Range("N98").Value = -4246
$ws.Range("H122").Value = 1346.4286
$ws.Range("J122").Value = 1250
$ws.Range("L122").Value = 3750
$ws.Range("N122").Value = -8650
$ws.Range("H137").Value = 963726.0600000001
$ws.Range("I137").Value = 27297.477
$ws.Range("J137").Value = 3148726
$ws.Range("K137").Value = 81892.431
$ws.Range("L137").Value = 9446178
$ws.Range("M137").Value = -79342.431
$ws.Range("N137").Value = -9451278
$ws.Range("H138").Value = 3924.943
$ws.Range("I138").Value = 1586.0667
$ws.Range("J138").Value = 4405.534
$ws.Range("K138").Value = 4758.2001
$ws.Range("L138").Value = 13216.602
$ws.Range("M138").Value = 381.7999
$ws.Range("N138").Value = -23496.602

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3776.5405
$ws.Range("I32").Value = 3449.742
$ws.Range("K32").Value = 3449.742
$ws.Range("M32").Value = -3162.742
$ws.Range("H61").Value = 3734.1428
$ws.Range("I61").Value = 3482.923
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 3482.923
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -3270.923
$ws.Range("N61").Value = -7424
$ws.Range("H110").Value = 10160.1
$ws.Range("I110").Value = 10949.5
$ws.Range("J110").Value = 3055.5
$ws.Range("K110").Value = 10949.5
$ws.Range("L110").Value = 3055.5
$ws.Range("M110").Value = -8904.5
$ws.Range("N110").Value = -7145.5
$ws.Range("H132").Value = 3273.7144
$ws.Range("I132").Value = 2776.1
$ws.Range("K132").Value = 8328.299999999999
$ws.Range("M132").Value = -5798.299999999999
$ws.Range("H136").Value = 3734.1428
$ws.Range("I136").Value = 3482.923
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 10448.769
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -7898.769
$ws.Range("N136").Value = -26100
$ws.Range("H139").Value = 76422.336
$ws.Range("J139").Value = 76422.336
$ws.Range("L139").Value = 76422.336
$ws.Range("N139").Value = -86702.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3797.1428
$ws.Range("I20").Value = 3660.4
$ws.Range("J20").Value = 4139
$ws.Range("K20").Value = 3660.4
$ws.Range("L20").Value = 4139
$ws.Range("M20").Value = -3413.4
$ws.Range("N20").Value = -4633
$ws.Range("H22").Value = 286.2
$ws.Range("I22").Value = 284
$ws.Range("K22").Value = 284
$ws.Range("M22").Value = -111
$ws.Range("H81").Value = 26010.6
$ws.Range("J81").Value = 26010.6
$ws.Range("L81").Value = 26010.6
$ws.Range("N81").Value = -28132.6
$ws.Range("H84").Value = 26010.6
$ws.Range("J84").Value = 26010.6
$ws.Range("L84").Value = 78031.79999999999
$ws.Range("N84").Value = -88639.79999999999
$ws.Range("H86").Value = 3962.9614
$ws.Range("I86").Value = 2638.4375
$ws.Range("J86").Value = 6082.2
$ws.Range("K86").Value = 2638.4375
$ws.Range("L86").Value = 6082.2
$ws.Range("M86").Value = -1515.4375
$ws.Range("N86").Value = -8328.200000000001
$ws.Range("H89").Value = 3962.9614
$ws.Range("I89").Value = 2638.4375
$ws.Range("J89").Value = 6082.2
$ws.Range("K89").Value = 13192.1875
$ws.Range("L89").Value = 30411
$ws.Range("M89").Value = -7576.1875
$ws.Range("N89").Value = -41643
$ws.Range("H94").Value = 1757.3077
$ws.Range("I94").Value = 1845.5
$ws.Range("J94").Value = 1463.3334
$ws.Range("K94").Value = 1845.5
$ws.Range("L94").Value = 1463.3334
$ws.Range("M94").Value = -1394.5
$ws.Range("N94").Value = -2365.3334
$ws.Range("H99").Value = 6945.769
$ws.Range("J99").Value = 10152.143
$ws.Range("L99").Value = 10152.143
$ws.Range("N99").Value = -13148.143
$ws.Range("H107").Value = 2624.5454
$ws.Range("I107").Value = 2176.4285
$ws.Range("K107").Value = 2176.4285
$ws.Range("M107").Value = -256.4285
$ws.Range("H140").Value = 71643.89999999999
$ws.Range("J140").Value = 71643.89999999999
$ws.Range("L140").Value = 71643.89999999999
$ws.Range("N140").Value = -82003.89999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 189.4762
$ws.Range("I7").Value = 107
$ws.Range("J7").Value = 323.5
$ws.Range("K7").Value = 107
$ws.Range("L7").Value = 323.5
$ws.Range("M7").Value = 6
$ws.Range("N7").Value = -549.5
$ws.Range("H94").Value = 1522.2858
$ws.Range("J94").Value = 499
$ws.Range("L94").Value = 499
$ws.Range("N94").Value = -1401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 90914860
$ws.Range("I3").Value = 142858080
$ws.Range("K3").Value = 428574240
$ws.Range("M3").Value = -428574128
$ws.Range("H68").Value = 1853360.6
$ws.Range("J68").Value = 2382649.2
$ws.Range("L68").Value = 7147947.600000001
$ws.Range("N68").Value = -7149569.600000001
$ws.Range("H71").Value = 1853360.6
$ws.Range("J71").Value = 2382649.2
$ws.Range("L71").Value = 21443842.8
$ws.Range("N71").Value = -21451954.8
$ws.Range("H107").Value = 27778584
$ws.Range("I107").Value = 33333840
$ws.Range("K107").Value = 100001520
$ws.Range("M107").Value = -99999600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6057.6
$ws.Range("J70").Value = 6310.857
$ws.Range("L70").Value = 6310.857
$ws.Range("N70").Value = -6850.857
$ws.Range("H73").Value = 6057.6
$ws.Range("J73").Value = 6310.857
$ws.Range("L73").Value = 6310.857
$ws.Range("N73").Value = -8182.857
$ws.Range("H97").Value = 6413.75
$ws.Range("I97").Value = 6258.5713
$ws.Range("K97").Value = 6258.5713
$ws.Range("M97").Value = -5762.5713
$ws.Range("H102").Value = 3256.8
$ws.Range("I102").Value = 3256.8
$ws.Range("K102").Value = 3256.8
$ws.Range("M102").Value = -1634.8
$ws.Range("H126").Value = 21426.611
$ws.Range("I126").Value = 38311.11
$ws.Range("K126").Value = 114933.33
$ws.Range("M126").Value = -112463.33
$ws.Range("H132").Value = 23403.875
$ws.Range("I132").Value = 25167.863
$ws.Range("K132").Value = 75503.58900000001
$ws.Range("M132").Value = -72973.58900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5308.1665
$ws.Range("J132").Value = 6690
$ws.Range("L132").Value = 20070
$ws.Range("N132").Value = -25130
$ws.Range("H136").Value = 3004
$ws.Range("I136").Value = 3214.8
$ws.Range("J136").Value = 1950
$ws.Range("K136").Value = 9644.400000000001
$ws.Range("L136").Value = 5850
$ws.Range("M136").Value = -7094.400000000001
$ws.Range("N136").Value = -10950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 1000
$ws.Range("I26").Value = 1000
$ws.Range("K26").Value = 1000
$ws.Range("M26").Value = -707
$ws.Range("H63").Value = 21999.285
$ws.Range("I63").Value = 19998
$ws.Range("J63").Value = 22332.834
$ws.Range("K63").Value = 19998
$ws.Range("L63").Value = 22332.834
$ws.Range("M63").Value = -19374
$ws.Range("N63").Value = -23580.834
$ws.Range("H66").Value = 21999.285
$ws.Range("I66").Value = 19998
$ws.Range("J66").Value = 22332.834
$ws.Range("K66").Value = 59994
$ws.Range("L66").Value = 66998.50199999999
$ws.Range("M66").Value = -56874
$ws.Range("N66").Value = -73238.50199999999
$ws.Range("H132").Value = 9196.200000000001
$ws.Range("I132").Value = 7345.55
$ws.Range("K132").Value = 22036.65
$ws.Range("M132").Value = -19506.65
